# Apply updated cryptocurrency price/volume data - Updated cryptos list on Fri Jun 23 08:07:43 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.915.88"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.879.00"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'242.73"
$ws.Range("E5").Value = "  -3.73%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.4913"
$ws.Range("E7").Value = "  -3.84%  "
$ws.Range("D8").Value = "'0.2932"
$ws.Range("E8").Value = "  -2.45%  "
$ws.Range("D9").Value = "'0.06640"
$ws.Range("E9").Value = "  -2.61%  "
$ws.Range("D10").Value = "1.875.73"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").Value = "'16.74"
$ws.Range("E11").Value = "  -3.62%  "
$ws.Range("D12").Value = "'0.07225"
$ws.Range("D13").Value = "'0.6676"
$ws.Range("E13").Value = "  -5.34%  "
$ws.Range("D14").Value = "'86.28"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").Value = "'4.880"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "29.906.70"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "'0.000007852"
$ws.Range("E17").Value = "  -4.04%  "
$ws.Range("D18").Value = "'0.9993"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("D20").Value = "2.119.83"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("D23").Value = "'5.830"
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("D24").Value = "'9.075"
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("D25").Value = "'150.05"
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("D26").Value = "'141.92"
$ws.Range("E26").Value = "  +4.85%  "
$ws.Range("D27").Value = "'17.04"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").Value = "'1.917"
$ws.Range("E28").Value = "  -4.43%  "
$ws.Range("D29").Value = "'1.386"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").Value = "'4.201"
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("D31").Value = "'0.08752"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").Value = "'3.968"
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("D33").Value = "'0.05044"
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("D34").Value = "'0.7126"
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").Value = "'1.112"
$ws.Range("E35").Value = "  -3.08%  "
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").Value = "'0.01812"
$ws.Range("E37").Value = "  +6.46%  "
$ws.Range("D38").Value = "'2.688"
$ws.Range("E38").Value = "  -4.65%  "
$ws.Range("D39").Value = "'2.171"
$ws.Range("E39").Value = "  -4.55%  "
$ws.Range("D40").Value = "'0.9299"
$ws.Range("E40").Value = "  -3.66%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'0.9992"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "'5.758"
$ws.Range("E42").Value = "  -6.89%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4232"
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("D44").Value = "'103.07"
$ws.Range("E44").Value = "  -2.53%  "
$ws.Range("D45").Value = "'7.415"
$ws.Range("E45").Value = "  -3.06%  "
$ws.Range("D46").Value = "'0.1269"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").Value = "'0.05666"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "'32.66"
$ws.Range("E48").Value = "  -2.26%  "
$ws.Range("D49").Value = "'0.3781"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").Value = "'8.289"
$ws.Range("E50").Value = "  -2.52%  "
$ws.Range("D51").Value = "'56.02"
$ws.Range("E51").Value = "  -1.52%  "
